$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.331.64'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '1.612.27'
$ws.Range('E3').Value = '  +0.36%  '

$ws.Range('D5').Value = "'213.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('E6').Value = '  -0.15%  '

$ws.Range('E7').Value = '  +0.29%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').Value = "'0.0616"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.33%  '

$ws.Range('D10').Value = "'18.53"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.49%  '

$ws.Range('D11').Value = "'0.0814"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.28%  '

$ws.Range('D12').Value = '1.833.98'

$ws.Range('D13').Value = '1.633.17'
$ws.Range('E13').Value = '  +1.66%  '

$ws.Range('E14').Value = '  +0.76%  '

$ws.Range('D15').Value = "'0.517"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.95%  '

$ws.Range('D16').Value = '26.325.51'
$ws.Range('E16').Value = '  +0.41%  '

$ws.Range('D17').Value = "'61.93"
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  +0.51%  '

$ws.Range('D20').Value = "'203.40"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.42%  '

$ws.Range('D21').Value = "'4.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.27%  '

$ws.Range('D22').Value = "'9.33"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.38%  '

$ws.Range('D23').Value = "'6.03"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.31%  '

$ws.Range('E24').Value = '  +8.64%  '

$ws.Range('D25').Value = "'144.32"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.29%  '

$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('E27').Value = '  -3.11%  '

$ws.Range('D28').Value = "'15.27"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.68%  '

$ws.Range('D29').Value = "'6.57"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.44%  '

$ws.Range('E30').Value = '  +3.87%  '

$ws.Range('D31').Value = "'1.18"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.24%  '

$ws.Range('E32').Value = '  +2.29%  '

$ws.Range('D33').Value = "'2.96"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.59%  '

$ws.Range('D34').Value = "'2.44"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.49%  '

$ws.Range('E35').Value = '  +0.49%  '

$ws.Range('D36').Value = '1.162.85'
$ws.Range('E36').Value = '  +4.98%  '

$ws.Range('D37').Value = "'0.0166"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.18%  '

$ws.Range('E38').Value = '  -0.12%  '

$ws.Range('D39').Value = "'0.798"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.46%  '

$ws.Range('E40').Value = '  -0.58%  '

$ws.Range('E41').Value = '  +0.64%  '

$ws.Range('D42').Value = "'0.788"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.67%  '

$ws.Range('E43').Value = '  +2.98%  '

$ws.Range('D44').Value = '1.747.50'
$ws.Range('E44').Value = '  +0.31%  '

$ws.Range('D45').Value = "'91.88"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.54%  '

$ws.Range('D46').Value = "'1.54"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.07%  '

$ws.Range('D47').Value = "'54.51"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.93%  '

$ws.Range('E48').Value = '  +0.06%  '

$ws.Range('D49').Value = '0.0₇0978'
$ws.Range('E49').Value = '  -14.14%  '

$ws.Range('D50').Value = "'0.407"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.65%  '

$ws.Range('E51').Value = '  -0.11%  '
